# Scheduled market-data refresh: update Leve profit columns (H:N) across job sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H53").Value = 199
$ws.Range("I53").Value = 206.71428
$ws.Range("J53").Value = 195.625
$ws.Range("K53").Value = 206.71428
$ws.Range("L53").Value = 195.625
$ws.Range("M53").Value = 430.28572
$ws.Range("N53").Value = -1469.625
$ws.Range("H138").Value = 2433.3096
$ws.Range("I138").Value = 1408.7872
$ws.Range("J138").Value = 3734.7297
$ws.Range("K138").Value = 4226.3616
$ws.Range("L138").Value = 11204.1891
$ws.Range("M138").Value = 913.6383999999998
$ws.Range("N138").Value = -21484.1891

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 753550.1
$ws.Range("I61").Value = 1103419.1
$ws.Range("J61").Value = 559178.5
$ws.Range("K61").Value = 1103419.1
$ws.Range("L61").Value = 559178.5
$ws.Range("M61").Value = -1103207.1
$ws.Range("N61").Value = -559602.5
$ws.Range("H63").Value = 3281.818
$ws.Range("I63").Value = 2900
$ws.Range("K63").Value = 2900
$ws.Range("M63").Value = -2214
$ws.Range("H66").Value = 3281.818
$ws.Range("I66").Value = 2900
$ws.Range("K66").Value = 14500
$ws.Range("M66").Value = -11068
$ws.Range("H122").Value = 3075.7073
$ws.Range("I122").Value = 2685.28
$ws.Range("J122").Value = 3685.75
$ws.Range("K122").Value = 8055.84
$ws.Range("L122").Value = 11057.25
$ws.Range("M122").Value = -5605.84
$ws.Range("N122").Value = -15957.25
$ws.Range("H136").Value = 753550.1
$ws.Range("I136").Value = 1103419.1
$ws.Range("J136").Value = 559178.5
$ws.Range("K136").Value = 3310257.3
$ws.Range("L136").Value = 1677535.5
$ws.Range("M136").Value = -3307707.3
$ws.Range("N136").Value = -1682635.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3120.6086
$ws.Range("I31").Value = 2266.423
$ws.Range("J31").Value = 4231.05
$ws.Range("K31").Value = 2266.423
$ws.Range("L31").Value = 4231.05
$ws.Range("M31").Value = -1971.423
$ws.Range("N31").Value = -4821.05
$ws.Range("H34").Value = 3120.6086
$ws.Range("I34").Value = 2266.423
$ws.Range("J34").Value = 4231.05
$ws.Range("K34").Value = 2266.423
$ws.Range("L34").Value = 4231.05
$ws.Range("M34").Value = -2064.423
$ws.Range("N34").Value = -4635.05
$ws.Range("H99").Value = 52212.3
$ws.Range("I99").Value = 84801
$ws.Range("J99").Value = 3329.25
$ws.Range("K99").Value = 84801
$ws.Range("L99").Value = 3329.25
$ws.Range("M99").Value = -83303
$ws.Range("N99").Value = -6325.25
$ws.Range("H107").Value = 312.8889
$ws.Range("I107").Value = 196.07692
$ws.Range("J107").Value = 616.6
$ws.Range("K107").Value = 196.07692
$ws.Range("L107").Value = 616.6
$ws.Range("M107").Value = 1723.92308
$ws.Range("N107").Value = -4456.6
$ws.Range("H122").Value = 2108.4546
$ws.Range("I122").Value = 1036
$ws.Range("J122").Value = 3002.1667
$ws.Range("K122").Value = 3108
$ws.Range("L122").Value = 9006.500100000001
$ws.Range("M122").Value = -658
$ws.Range("N122").Value = -13906.5001
$ws.Range("H126").Value = 52212.3
$ws.Range("I126").Value = 84801
$ws.Range("J126").Value = 3329.25
$ws.Range("K126").Value = 254403
$ws.Range("L126").Value = 9987.75
$ws.Range("M126").Value = -251933
$ws.Range("N126").Value = -14927.75
$ws.Range("H132").Value = 2486.6453
$ws.Range("I132").Value = 1202.75
$ws.Range("J132").Value = 3856.1333
$ws.Range("K132").Value = 3608.25
$ws.Range("L132").Value = 11568.3999
$ws.Range("M132").Value = -1078.25
$ws.Range("N132").Value = -16628.3999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H88").Value = 5625
$ws.Range("J88").Value = 5625
$ws.Range("L88").Value = 16875
$ws.Range("N88").Value = -17731
$ws.Range("H91").Value = 5625
$ws.Range("J91").Value = 5625
$ws.Range("L91").Value = 16875
$ws.Range("N91").Value = -19839
$ws.Range("H122").Value = 852.069
$ws.Range("I122").Value = 412.07693
$ws.Range("J122").Value = 1209.5625
$ws.Range("K122").Value = 3708.69237
$ws.Range("L122").Value = 10886.0625
$ws.Range("M122").Value = -1258.69237
$ws.Range("N122").Value = -15786.0625
$ws.Range("H132").Value = 6997.5884
$ws.Range("I132").Value = 8100
$ws.Range("J132").Value = 6761.357
$ws.Range("K132").Value = 72900
$ws.Range("L132").Value = 60852.213
$ws.Range("M132").Value = -70370
$ws.Range("N132").Value = -65912.213

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 1544.9546
$ws.Range("I122").Value = 1269.9412
$ws.Range("J122").Value = 2480
$ws.Range("K122").Value = 3809.8236
$ws.Range("L122").Value = 7440
$ws.Range("M122").Value = -1359.8236
$ws.Range("N122").Value = -12340
$ws.Range("H126").Value = 5618.625
$ws.Range("I126").Value = 5040
$ws.Range("J126").Value = 6583
$ws.Range("K126").Value = 15120
$ws.Range("L126").Value = 19749
$ws.Range("M126").Value = -12650
$ws.Range("N126").Value = -24689

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H5").Value = 5943
$ws.Range("I5").Value = 3914.5
$ws.Range("K5").Value = 3914.5
$ws.Range("M5").Value = -3801.5
$ws.Range("H7").Value = 2233.3872
$ws.Range("I7").Value = 1905.55
$ws.Range("J7").Value = 2829.4546
$ws.Range("K7").Value = 1905.55
$ws.Range("L7").Value = 2829.4546
$ws.Range("M7").Value = -1793.55
$ws.Range("N7").Value = -3053.4546
$ws.Range("H122").Value = 2218.4375
$ws.Range("I122").Value = 2308.182
$ws.Range("J122").Value = 2021
$ws.Range("K122").Value = 6924.545999999999
$ws.Range("L122").Value = 6063
$ws.Range("M122").Value = -4474.545999999999
$ws.Range("N122").Value = -10963
$ws.Range("H126").Value = 2233.3872
$ws.Range("I126").Value = 1905.55
$ws.Range("J126").Value = 2829.4546
$ws.Range("K126").Value = 5716.65
$ws.Range("L126").Value = 8488.363799999999
$ws.Range("M126").Value = -3246.65
$ws.Range("N126").Value = -13428.3638

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H21").Value = 10002923
$ws.Range("I21").Value = 16668372
$ws.Range("J21").Value = 4750
$ws.Range("K21").Value = 16668372
$ws.Range("L21").Value = 4750
$ws.Range("M21").Value = -16668137
$ws.Range("N21").Value = -5220
$ws.Range("H24").Value = 2005
$ws.Range("J24").Value = 2005
$ws.Range("L24").Value = 2005
$ws.Range("N24").Value = -2465
$ws.Range("H28").Value = 4943.4
$ws.Range("I28").Value = 3517
$ws.Range("J28").Value = 5300
$ws.Range("K28").Value = 3517
$ws.Range("L28").Value = 5300
$ws.Range("M28").Value = -3169
$ws.Range("N28").Value = -5996
$ws.Range("H35").Value = 10002923
$ws.Range("I35").Value = 16668372
$ws.Range("J35").Value = 4750
$ws.Range("K35").Value = 16668372
$ws.Range("L35").Value = 4750
$ws.Range("M35").Value = -16668082
$ws.Range("N35").Value = -5330
$ws.Range("H107").Value = 1052.8
$ws.Range("I107").Value = 1020.6667
$ws.Range("J107").Value = 1101
$ws.Range("K107").Value = 3062.0001
$ws.Range("L107").Value = 3303
$ws.Range("M107").Value = -1142.0001
$ws.Range("N107").Value = -7143
$ws.Range("H122").Value = 62502780
$ws.Range("I122").Value = 83336040
$ws.Range("J122").Value = 3001.25
$ws.Range("K122").Value = 250008120
$ws.Range("L122").Value = 9003.75
$ws.Range("M122").Value = -250005670
$ws.Range("N122").Value = -13903.75
$ws.Range("H126").Value = 974.21875
$ws.Range("I126").Value = 573.0769
$ws.Range("J126").Value = 2712.5
$ws.Range("K126").Value = 1719.2307
$ws.Range("L126").Value = 8137.5
$ws.Range("M126").Value = 750.7692999999999
$ws.Range("N126").Value = -13077.5
$ws.Range("H132").Value = 2384.1538
$ws.Range("I132").Value = 1185.0714
$ws.Range("J132").Value = 3783.0833
$ws.Range("K132").Value = 3555.2142
$ws.Range("L132").Value = 11349.2499
$ws.Range("M132").Value = -1025.2142
$ws.Range("N132").Value = -16409.2499

